$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D column holds numeric-looking price strings (e.g. "27.939.63", "329.90", "0.07843").
# Excel auto-coerces such literals to floating point numbers when assigned via .Value,
# which both changes the cell type and loses exact formatting (trailing zeros, sci notation).
# Forcing NumberFormat to Text ("@") before assignment keeps the literal string intact,
# matching how these price cells are stored as inline/shared strings in the workbook.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.939.63"
$ws.Range("E2").Value = "  -5.41%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.822.08"
$ws.Range("E3").Value = "  -4.45%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.74%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "329.90"
$ws.Range("E5").Value = "  -2.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  -0.66%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4629"
$ws.Range("E7").Value = "  -2.81%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3844"
$ws.Range("E8").Value = "  -3.88%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.84"
$ws.Range("E9").Value = "  -3.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07843"
$ws.Range("E10").Value = "  -2.41%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9601"
$ws.Range("E11").Value = "  -3.08%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.84"
$ws.Range("E12").Value = "  -6.25%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.785.30"
$ws.Range("E13").Value = "  -6.80%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.641"
$ws.Range("E14").Value = "  -4.71%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.853"
$ws.Range("E15").Value = "  -3.61%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06845"
$ws.Range("E16").Value = "  +0.11%  "
$ws.Range("E17").Value = "  -0.73%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "86.66"
$ws.Range("E18").Value = "  -2.75%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000009937"
$ws.Range("E19").Value = "  -2.73%  "
$ws.Range("E20").Value = "  -3.91%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.003"
$ws.Range("E21").Value = "  -0.40%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.961.16"
$ws.Range("E22").Value = "  -5.34%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.317"
$ws.Range("E23").Value = "  -3.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.94"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.099"
$ws.Range("E25").Value = "  -2.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.058.64"
$ws.Range("E26").Value = "  -3.98%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "152.31"
$ws.Range("E27").Value = "  -3.01%  "
$ws.Range("E28").Value = "  -2.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.693"
$ws.Range("E29").Value = "  -12.85%  "
$ws.Range("E30").Value = "  -4.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "116.52"
$ws.Range("E31").Value = "  -2.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9366"
$ws.Range("E32").Value = "  -5.81%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09259"
$ws.Range("E33").Value = "  -2.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.268"
$ws.Range("E34").Value = "  -3.94%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.418"
$ws.Range("E35").Value = "  -3.62%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.312"
$ws.Range("E36").Value = "  -5.55%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05933"
$ws.Range("E37").Value = "  -8.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02148"
$ws.Range("E38").Value = "  -4.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.147"
$ws.Range("E39").Value = "  -3.70%  "
$ws.Range("E40").Value = "  -0.63%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.570"
$ws.Range("E41").Value = "  -2.26%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5575"
$ws.Range("E42").Value = "  -4.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "9.898"
$ws.Range("E43").Value = "  -5.74%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1769"
$ws.Range("E44").Value = "  -2.93%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.222"
$ws.Range("E45").Value = "  -3.89%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.228"
$ws.Range("E46").Value = "  -9.12%  "
$ws.Range("E47").Value = "  -5.09%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5244"
$ws.Range("E48").Value = "  -4.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07005"
$ws.Range("E49").Value = "  -5.66%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.822"
$ws.Range("E50").Value = "  -6.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "112.60"
$ws.Range("E51").Value = "  -3.03%  "
